$wb = $excel.ActiveWorkbook

# Insert a new worksheet "Sheet2" right after the existing "Sheet1"
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate the new sheet with the dynamic search data
$ws2.Range("A1").Value = "Searchitem"
$ws2.Range("A2").Value = "65 inch tv Samsung"
$ws2.Range("A3").Value = "65 inch tv lg"

# Widen column A to fit the new content
$ws2.Columns.Item(1).ColumnWidth = 18

# Match the saved selection/view state for the new sheet
$ws2.Range("F5").Select() | Out-Null

# Make Sheet2 the active (visible) tab
$ws2.Activate() | Out-Null
